# Update column F (dSF) values for several rows based on repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    19 = 3
    25 = 4
    28 = -1
    30 = 0
    46 = 2
    48 = 1
    55 = 1
    71 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
